# MAJ collecteurs (prix + chemins)
# Renames the image-path folder "collecteur/" -> "collecteurs/" for the
# collector rows whose picture moved (rows 2 & 3 keep their original
# "collecteur/" folder, matching the source diff), then restores the
# selection / scroll position the author left the sheet in.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A ("path") updates ---------------------------------------
# Rows 2 and 3 are unchanged (still "collecteur/...").
$ws.Range("A4").Value  = "collecteurs/collecteur-4-rouge.png"
$ws.Range("A5").Value  = "collecteurs/collecteur-5-rouge.png"
$ws.Range("A6").Value  = "collecteurs/collecteur-6-rouge.png"
$ws.Range("A7").Value  = "collecteurs/collecteur-2-bleu.png"
$ws.Range("A8").Value  = "collecteurs/collecteur-3-bleu.png"
$ws.Range("A9").Value  = "collecteurs/collecteur-4-bleu.png"
$ws.Range("A10").Value = "collecteurs/collecteur-5-bleu.png"
$ws.Range("A11").Value = "collecteurs/collecteur-6-bleu.png"
$ws.Range("A19").Value = "collecteurs/collecteur-2-bleu.png"
$ws.Range("A20").Value = "collecteurs/collecteur-2-rouge.png"
$ws.Range("A21").Value = "collecteurs/collecteur-5-rouge.png"
$ws.Range("A22").Value = "collecteurs/collecteur-5-bleu.png"

# --- Column A width (bestFit grew slightly because the longest path
#     string grew by one character) --------------------------------
$ws.Columns.Item(1).ColumnWidth = 29.44140625

# --- Restore cursor / view position left by the author --------------
$ws.Range("A22").Select()
